# Temp Screening, Visitor Purpose Related Questions, Updates Screening Questionnaire
# Reassign three sub-tasks on "Sprint 1 (17Oct-31Oct)" from Aloysius to
# Christopher (H9, H11, H12). The dependent per-person marker formulas in
# columns I:N recalc automatically. Also update the sheet's view/selection
# state to match (active cell H10, scrolled so row 4 is at the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1 (17Oct-31Oct)")
$ws.Activate()

foreach ($r in 9, 11, 12) {
    $cell = $ws.Cells.Item($r, 8)   # column H = "Assigned To"
    $cell.Value = "Christopher"
    # These rows sit away from the header, but Excel applied the
    # "just-below-header" border variant (no top border) when the value was
    # re-entered, so drop the top border to match.
    $cell.Borders.Item(8).LineStyle = 0
}

# Move the view so row 4 is the first visible row, and land the selection
# on H10.
$ws.Application.Goto($ws.Range("A4"), $false)
$ws.Range("H10").Select()
